# Update the crypto symbol list (price / volume refresh + a couple of
# row insertions/shifts for the exchange-token section) as produced by the
# "Updated symbol list" GitHub Actions run.
#
# Price (column D) and Volume(1h) (column E) values are text that merely
# look numeric (e.g. "304.70", "-0.69%"); they are written with a leading
# apostrophe so Excel stores them as text instead of silently coercing them
# to Number/Percentage values (which would also lose the original,
# significant trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.70"
$ws.Range("E2").Value = "'-0.69%"
$ws.Range("D3").Value = "'35.91"
$ws.Range("E3").Value = "'-1.30%"
$ws.Range("D4").Value = "'5.022"
$ws.Range("E4").Value = "'-1.38%"
$ws.Range("D5").Value = "'0.08057"
$ws.Range("E5").Value = "'-0.27%"
$ws.Range("D6").Value = "'1.880"
$ws.Range("E6").Value = "'-4.45%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.153"
$ws.Range("E7").Value = "'0.90%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.842"
$ws.Range("E8").Value = "'1.10%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9315"
$ws.Range("E9").Value = "'-0.39%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1306"
$ws.Range("E10").Value = "'-9.20%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1898"
$ws.Range("E11").Value = "'-1.42%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09183"
$ws.Range("E12").Value = "'-0.10%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03517"
$ws.Range("E13").Value = "'-0.93%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09899"
$ws.Range("E14").Value = "'1.06%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001423"
$ws.Range("E15").Value = "'-1.14%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006329"
$ws.Range("E16").Value = "'8.67%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.621"
$ws.Range("E17").Value = "'2.71%"
$ws.Range("D18").Value = "'3.220"
$ws.Range("E18").Value = "'9.83%"
$ws.Range("E19").Value = "'0.49%"
$ws.Range("E20").Value = "'2.41%"
$ws.Range("D21").Value = "'5.219"
$ws.Range("E21").Value = "'5.17%"
$ws.Range("D22").Value = "'0.2535"
$ws.Range("E22").Value = "'5.20%"
$ws.Range("D23").Value = "'0.04421"
$ws.Range("E23").Value = "'-2.10%"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'1.76%"
$ws.Range("D25").Value = "'0.004704"
$ws.Range("E25").Value = "'-3.00%"
$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'4.83%"
$ws.Range("E27").Value = "'-29.70%"
$ws.Range("D39").Value = "'0.01950"
$ws.Range("E39").Value = "'-1.21%"
$ws.Range("D40").Value = "'0.05160"
$ws.Range("E40").Value = "'6.37%"
$ws.Range("E41").Value = "'0.13%"
$ws.Range("D42").Value = "'0.01017"
$ws.Range("E42").Value = "'-9.46%"
$ws.Range("D43").Value = "'0.1370"
$ws.Range("E43").Value = "'0.52%"
$ws.Range("D44").Value = "'0.002174"
$ws.Range("E44").Value = "'4.36%"
$ws.Range("D45").Value = "'0.01077"
$ws.Range("E45").Value = "'10.40%"
$ws.Range("D46").Value = "'0.00006352"
$ws.Range("E46").Value = "'-0.09%"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("D48").Value = "'63.57"
$ws.Range("E48").Value = "'-1.69%"
$ws.Range("D49").Value = "'0.001662"
$ws.Range("E49").Value = "'39.29%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'-0.04%"
